$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4276176
$ws.Range("I33").Value = 6536628
$ws.Range("K33").Value = 6536628
$ws.Range("M33").Value = -6536399

$ws.Range("H40").Value = 1954.6
$ws.Range("I40").Value = 1863.2
$ws.Range("J40").Value = 2009.44
$ws.Range("K40").Value = 1863.2
$ws.Range("L40").Value = 2009.44
$ws.Range("M40").Value = -1688.2
$ws.Range("N40").Value = -2359.44

$ws.Range("H86").Value = 4103.625
$ws.Range("I86").Value = 4225.5
$ws.Range("K86").Value = 4225.5
$ws.Range("M86").Value = -3102.5

$ws.Range("H89").Value = 4103.625
$ws.Range("I89").Value = 4225.5
$ws.Range("K89").Value = 21127.5
$ws.Range("M89").Value = -15511.5

$ws.Range("H113").Value = 4762.913
$ws.Range("I113").Value = 3797.4375
$ws.Range("K113").Value = 3797.4375
$ws.Range("M113").Value = -543.4375

$ws.Range("H137").Value = 9875.9
$ws.Range("I137").Value = 4533.0454
$ws.Range("J137").Value = 14073.857
$ws.Range("K137").Value = 13599.1362
$ws.Range("L137").Value = 42221.571
$ws.Range("M137").Value = -11049.1362
$ws.Range("N137").Value = -47321.571

$ws.Range("H141").Value = 3006
$ws.Range("I141").Value = 2889.8
$ws.Range("J141").Value = 3151.25
$ws.Range("K141").Value = 8669.400000000001
$ws.Range("L141").Value = 9453.75
$ws.Range("M141").Value = -3489.400000000001
$ws.Range("N141").Value = -19813.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4181.2
$ws.Range("I32").Value = 3718.6155
$ws.Range("K32").Value = 3718.6155
$ws.Range("M32").Value = -3431.6155

$ws.Range("H45").Value = 16643.555
$ws.Range("I45").Value = 16643.555
$ws.Range("K45").Value = 16643.555
$ws.Range("M45").Value = -16266.555

$ws.Range("H61").Value = 11329.607
$ws.Range("I61").Value = 9815.4
$ws.Range("K61").Value = 9815.4
$ws.Range("M61").Value = -9603.4

$ws.Range("H74").Value = 10602.88
$ws.Range("I74").Value = 11148.869
$ws.Range("K74").Value = 11148.869
$ws.Range("M74").Value = -10274.869

$ws.Range("H77").Value = 10602.88
$ws.Range("I77").Value = 11148.869
$ws.Range("K77").Value = 55744.345
$ws.Range("M77").Value = -51376.345

$ws.Range("H105").Value = 57122.668
$ws.Range("J105").Value = 57122.668
$ws.Range("L105").Value = 57122.668
$ws.Range("N105").Value = -64110.668

$ws.Range("H136").Value = 11329.607
$ws.Range("I136").Value = 9815.4
$ws.Range("K136").Value = 29446.2
$ws.Range("M136").Value = -26896.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 72858.38
$ws.Range("J132").Value = 72858.38
$ws.Range("L132").Value = 72858.38
$ws.Range("N132").Value = -82978.38

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2631.9565
$ws.Range("I16").Value = 2846.25
$ws.Range("K16").Value = 2846.25
$ws.Range("M16").Value = -2559.25

$ws.Range("H31").Value = 3990.2222
$ws.Range("I31").Value = 1860.3334
$ws.Range("J31").Value = 5055.1665
$ws.Range("K31").Value = 1860.3334
$ws.Range("L31").Value = 5055.1665
$ws.Range("M31").Value = -1565.3334
$ws.Range("N31").Value = -5645.1665

$ws.Range("H33").Value = 1998.8
$ws.Range("J33").Value = 1998.5
$ws.Range("L33").Value = 1998.5
$ws.Range("N33").Value = -2756.5

$ws.Range("H34").Value = 3990.2222
$ws.Range("I34").Value = 1860.3334
$ws.Range("J34").Value = 5055.1665
$ws.Range("K34").Value = 1860.3334
$ws.Range("L34").Value = 5055.1665
$ws.Range("M34").Value = -1658.3334
$ws.Range("N34").Value = -5459.1665

$ws.Range("H107").Value = 1976.8572
$ws.Range("I107").Value = 1972
$ws.Range("J107").Value = 1983.3334
$ws.Range("K107").Value = 1972
$ws.Range("L107").Value = 1983.3334
$ws.Range("M107").Value = -52
$ws.Range("N107").Value = -5823.3334

$ws.Range("H113").Value = 2631.9565
$ws.Range("I113").Value = 2846.25
$ws.Range("K113").Value = 2846.25
$ws.Range("M113").Value = -676.25

$ws.Range("H132").Value = 21012.062
$ws.Range("I132").Value = 12347.105
$ws.Range("J132").Value = 33676.23
$ws.Range("K132").Value = 37041.315
$ws.Range("L132").Value = 101028.69
$ws.Range("M132").Value = -34511.315
$ws.Range("N132").Value = -106088.69

$ws.Range("H134").Value = 5228.3667
$ws.Range("I134").Value = 3274.28
$ws.Range("K134").Value = 9822.84
$ws.Range("M134").Value = -7287.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4557.5
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 35675.555
$ws.Range("J57").Value = 53328.332
$ws.Range("L57").Value = 53328.332
$ws.Range("N57").Value = -54968.332

$ws.Range("H80").Value = 1239.0834
$ws.Range("I80").Value = 1240.6
$ws.Range("J80").Value = 1238
$ws.Range("K80").Value = 1240.6
$ws.Range("L80").Value = 1238
$ws.Range("M80").Value = -242.5999999999999
$ws.Range("N80").Value = -3234

$ws.Range("H83").Value = 1239.0834
$ws.Range("I83").Value = 1240.6
$ws.Range("J83").Value = 1238
$ws.Range("K83").Value = 6203
$ws.Range("L83").Value = 6190
$ws.Range("M83").Value = -1211
$ws.Range("N83").Value = -16174

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = ""

$ws.Range("H113").Value = 79115.96000000001
$ws.Range("I113").Value = 102142.7
$ws.Range("K113").Value = 102142.7
$ws.Range("M113").Value = -99972.7

$ws.Range("H126").Value = 4649.778
$ws.Range("I126").Value = 4207.6665
$ws.Range("J126").Value = 5534
$ws.Range("K126").Value = 12622.9995
$ws.Range("L126").Value = 16602
$ws.Range("M126").Value = -10152.9995
$ws.Range("N126").Value = -21542

$ws.Range("H132").Value = 12527.274
$ws.Range("I132").Value = 10488.214
$ws.Range("J132").Value = 16809.3
$ws.Range("K132").Value = 31464.642
$ws.Range("L132").Value = 50427.89999999999
$ws.Range("M132").Value = -28934.642
$ws.Range("N132").Value = -55487.89999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2338.913
$ws.Range("J22").Value = 2439.5881
$ws.Range("L22").Value = 2439.5881
$ws.Range("N22").Value = -3029.5881

$ws.Range("H27").Value = 2338.913
$ws.Range("J27").Value = 2439.5881
$ws.Range("L27").Value = 2439.5881
$ws.Range("N27").Value = -2653.5881

$ws.Range("H46").Value = 1532.5
$ws.Range("I46").Value = 951.9524
$ws.Range("J46").Value = 2249.647
$ws.Range("K46").Value = 951.9524
$ws.Range("L46").Value = 2249.647
$ws.Range("M46").Value = -763.9524
$ws.Range("N46").Value = -2625.647

$ws.Range("H122").Value = 3919.5
$ws.Range("I122").Value = 4522.3335
$ws.Range("K122").Value = 13567.0005
$ws.Range("M122").Value = -11117.0005

$ws.Range("H136").Value = 3365.3215
$ws.Range("I136").Value = 3010.0715
$ws.Range("K136").Value = 9030.2145
$ws.Range("M136").Value = -6480.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = ""
$ws.Range("N33").Value = ""

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = ""
$ws.Range("N36").Value = ""

$ws.Range("H37").Value = 83330.664
$ws.Range("J37").Value = 99996.5
$ws.Range("L37").Value = 99996.5
$ws.Range("N37").Value = -100402.5

$ws.Range("H103").Value = 16900.25
$ws.Range("J103").Value = 16900.25
$ws.Range("L103").Value = 16900.25
$ws.Range("N103").Value = -19244.25

$ws.Range("H136").Value = 2988720.5
$ws.Range("I136").Value = 4351537
$ws.Range("J136").Value = 3503.476
$ws.Range("K136").Value = 13054611
$ws.Range("L136").Value = 10510.428
$ws.Range("M136").Value = -13052061
$ws.Range("N136").Value = -15610.428
